# Auto commit at 2026-01-05  8:00:19.81
# Updates the monthly metric figures on the "Metrics" sheet (which cascade
# via formulas into the "today" sheet), and refreshes the remembered cell
# selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Sheets.Item("Metrics")
$wsToday   = $wb.Sheets.Item("today")

# --- Update the source figures on the Metrics sheet -----------------------
$wsMetrics.Range("B2").Value  = 72157.510000000009
$wsMetrics.Range("B3").Value  = 55508
$wsMetrics.Range("B4").Value  = 18397.52
$wsMetrics.Range("B5").Value  = 2900
$wsMetrics.Range("B6").Value  = 5656411.5900000008
$wsMetrics.Range("B7").Value  = 4786614.42
$wsMetrics.Range("B8").Value  = 1669294.25
$wsMetrics.Range("B9").Value  = 221079
$wsMetrics.Range("B10").Value = 34121792.579999998
$wsMetrics.Range("B11").Value = 32061889.579999998
$wsMetrics.Range("B12").Value = 11951016.289999999
$wsMetrics.Range("B13").Value = 1318709

# --- Restore the remembered selection on the Metrics sheet -----------------
# (Activate it only long enough to move its cached selection; the "today"
# sheet is re-activated afterwards so it remains the workbook's active tab.)
[void]$wsMetrics.Activate()
[void]$wsMetrics.Range("D16").Select()

# --- Restore the remembered selection on the "today" sheet (stays active) --
[void]$wsToday.Activate()
[void]$wsToday.Range("F8").Select()
